$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.150.87"
$ws.Range("E2").Value = "  +6.22%  "
$ws.Range("D3").Value = "3.538.91"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "417.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +6.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.788"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +20.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.86"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000274"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +28.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +10.00%  "
$ws.Range("D14").Value = "4.071.39"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.140"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.57"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "3.532.20"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").Value = "65.996.54"
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "449.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.27"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.20"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.96"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.39"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.57"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.34"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0506"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").Value = "0.0₃0741"
$ws.Range("E37").Value = "  +43.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.148"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.15%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  -5.06%  "
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("E48").Value = "  +6.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("E50").Value = "  +8.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.81%  "
